$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New poll data rows (ifop poll, 11/8)
$data = @(
    @{ A=28; B=2021; C=10; D=11; E=4; F="ifop"; G="online"; H="included"; I=1368; J=0.5;  K=0.5; L=8.5; M=2;   N=2.5; O=7; P=5; Q=25;              T=13; U=0.5; V=2.5; W=16;   X=17 },
    @{ A=28; B=2021; C=10; D=11; E=4; F="ifop"; G="online"; H="included"; I=1368; J=1;    K=0.5; L=8;   M=2.5; N=3;   O=7; P=6; Q=25;        S=10;       U=0.5; V=3;   W=17;   X=16.5 },
    @{ A=28; B=2021; C=10; D=11; E=4; F="ifop"; G="online"; H="included"; I=1368; J=1;    K=0.5; L=8;   M=2;   N=3;   O=7; P=6; Q=26; R=9;        U=1;   V=3.5; W=16.5; X=16.5 }
)

$startRow = 96
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}

# Update view state to match author's final selection/scroll position
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 88
$win.ScrollColumn = 4
$ws.Range("U99").Select()

$wb.Save()
